$wb = $excel.ActiveWorkbook

# --- Update the "random string" test data on the DemoWebShop sheet (row 2) ---
$ws = $wb.Worksheets.Item("DemoWebShop")

$ws.Range("C2").Value = "fnnxoWnO@gmail.com"
$ws.Range("F2").Value = "tUzBQ"
$ws.Range("G2").Value = "msSGx"

# --- Switch the active sheet / selection (framework change) ---
# Previously "OrangeHRM" was the selected tab; now "DemoWebShop" is active,
# with its own selection moved to C13.
$ws.Activate()
$ws.Range("C13").Select()
